$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row
# (rows 2-308). Bump it from 2023-09-20 (45189) to 2023-09-21 (45190) for
# every row, matching the bulk-update commit.
$ws.Range("C2:C308").Value = 45190
